$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Touch the format first so Excel stores these literals as text (matching the
# original inline-string cells) instead of re-parsing them as numbers.
$priceRange = $ws.Range("D2:E51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "26.088.40"
$ws.Range("D3").Value = "1.638.04"
$ws.Range("E3").Value = "  -1.92%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "213.92"
$ws.Range("E5").Value = "  +1.35%  "
$ws.Range("D6").Value = "0.5251"
$ws.Range("E6").Value = "  -0.23%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "0.2604"
$ws.Range("E8").Value = "  -1.20%  "
$ws.Range("D9").Value = "0.06310"
$ws.Range("E9").Value = "  -0.07%  "
$ws.Range("D10").Value = "20.72"
$ws.Range("E10").Value = "  -2.40%  "
$ws.Range("D11").Value = "0.07667"
$ws.Range("E11").Value = "  +1.29%  "
$ws.Range("D12").Value = "1.634.80"
$ws.Range("E12").Value = "  -2.79%  "
$ws.Range("D13").Value = "4.424"
$ws.Range("E13").Value = "  -0.49%  "
$ws.Range("D14").Value = "1.858.88"
$ws.Range("E14").Value = "  -2.38%  "
$ws.Range("D15").Value = "0.5512"
$ws.Range("E15").Value = "  -1.27%  "
$ws.Range("D16").Value = "0.0₅8129"
$ws.Range("E16").Value = "  +2.05%  "
$ws.Range("D17").Value = "65.19"
$ws.Range("E17").Value = "  -2.71%  "
$ws.Range("D18").Value = "26.073.12"
$ws.Range("E18").Value = "  -0.41%  "
$ws.Range("E19").Value = "  -0.08%  "
$ws.Range("D20").Value = "4.700"
$ws.Range("E20").Value = "  -1.04%  "
$ws.Range("D21").Value = "188.94"
$ws.Range("E21").Value = "  +0.94%  "
$ws.Range("D22").Value = "10.15"
$ws.Range("E22").Value = "  -2.31%  "
$ws.Range("D23").Value = "6.165"
$ws.Range("E23").Value = "  -0.36%  "
$ws.Range("D24").Value = "1.002"
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("D25").Value = "146.31"
$ws.Range("D26").Value = "0.1218"
$ws.Range("E26").Value = "  -2.66%  "
$ws.Range("D28").Value = "15.86"
$ws.Range("E28").Value = "  -0.92%  "
$ws.Range("E29").Value = "  +3.94%  "
$ws.Range("D30").Value = "0.06000"
$ws.Range("E30").Value = "  -4.48%  "
$ws.Range("D31").Value = "1.258"
$ws.Range("E31").Value = "  -1.81%  "
$ws.Range("D32").Value = "3.446"
$ws.Range("E32").Value = "  -1.85%  "
$ws.Range("D33").Value = "3.414"
$ws.Range("E33").Value = "  -0.12%  "
$ws.Range("D34").Value = "1.643"
$ws.Range("E34").Value = "  +0.81%  "
$ws.Range("D35").Value = "0.9883"
$ws.Range("E36").Value = "  -0.43%  "
$ws.Range("D37").Value = "2.762"
$ws.Range("E37").Value = "  +0.87%  "
$ws.Range("D38").Value = "0.5736"
$ws.Range("E38").Value = "  -5.38%  "
$ws.Range("D39").Value = "0.01618"
$ws.Range("E39").Value = "  +0.33%  "
$ws.Range("D40").Value = "0.8533"
$ws.Range("E40").Value = "  -2.66%  "
$ws.Range("D41").Value = "1.041.15"
$ws.Range("E41").Value = "  -5.69%  "
$ws.Range("D42").Value = "1.001"
$ws.Range("E42").Value = "  -0.25%  "
$ws.Range("D43").Value = "5.692"
$ws.Range("E43").Value = "  -6.98%  "
$ws.Range("D45").Value = "1.786.23"
$ws.Range("E45").Value = "  -2.00%  "
$ws.Range("E46").Value = "  -5.06%  "
$ws.Range("D47").Value = "55.42"
$ws.Range("E47").Value = "  -0.13%  "
$ws.Range("D48").Value = "0.9983"
$ws.Range("E48").Value = "  -0.13%  "
$ws.Range("D49").Value = "8.040"
$ws.Range("E49").Value = "  +0.08%  "
$ws.Range("D50").Value = "0.05171"
$ws.Range("E50").Value = "  -1.18%  "
$ws.Range("D51").Value = "0.4222"
$ws.Range("E51").Value = "  -0.60%  "

# Restore the default (unstyled) look so only the values themselves changed.
$priceRange.Style = "Normal"

